$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 11-14 - Posts 10-13 move from "In Queue" to "Edited"
$ws.Range("D11").Value = "Edited"
$ws.Range("D12").Value = "Edited"
$ws.Range("D13").Value = "Edited"
$ws.Range("D14").Value = "Edited"

# Row 10 - Post 9 "Beach Rock": add the image file name and mark as Posted
$ws.Range("C10").Value = "BeachRock.jpg"
$ws.Range("D10").Value = "Posted"

# Update the selected cell on the sheet
$ws.Range("J12").Select()
